# ===========================================================================
# Add "2022-Q4" quarterly sheet to the 09999-网易公司 (NetEase HK) workbook,
# positioned right after "总计" and before "2022-Q1", and record its summary
# row on the "总计" sheet.
# ===========================================================================

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the cell to be stored as text (keeps things like leading zeros
    # in fund codes, e.g. "012805", and decimal-looking figures such as
    # "246.27" from being silently re-typed as numbers).
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new first data row for
#    2022-Q4 and push the existing quarters down by one row.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Grow the table by one row, copying the formatting of the last existing
# data row (row 5, "2021-Q2") down onto the new row 6.
$wsTotal.Range("A5:D5").Copy($wsTotal.Range("A6:D6"))

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 27
$wsTotal.Range("D2").Value = 48.54

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q1"
$wsTotal.Range("C3").Value = 18
$wsTotal.Range("D3").Value = 25.87

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q4"
$wsTotal.Range("C4").Value = 11
$wsTotal.Range("D4").Value = 13.12

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2021-Q3"
$wsTotal.Range("C5").Value = 1
$wsTotal.Range("D5").Value = 2.96

$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").Value = "2021-Q2"
$wsTotal.Range("C6").Value = 1
$wsTotal.Range("D6").Value = 2.55

# ---------------------------------------------------------------------------
# 2) Create the new "2022-Q4" sheet. Cloning the existing "2022-Q1" sheet
#    gives us the right header row / column styling (bold + bordered header,
#    bold + bordered index column) for free; it is inserted immediately
#    before "2022-Q1", matching the desired tab order.
# ---------------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsQ1.Copy($wsQ1)
$wsQ4 = $wb.Worksheets.Item("2022-Q1 (2)")
$wsQ4.Name = "2022-Q4"

# The template has 19 rows (header + 18 funds); the new sheet needs 28 rows
# (header + 27 funds), so extend the table, copying the last template row's
# formatting down for each new row.
for ($r = 20; $r -le 28; $r++) {
    $wsQ4.Range("A19:H19").Copy($wsQ4.Range("A$r`:H$r"))
}


# ---------------------------------------------------------------------------
# 3) Populate the "2022-Q4" sheet's 27 fund-holding rows.
# ---------------------------------------------------------------------------
$wsQ4.Range("A2").Value = 0
Set-TextValue $wsQ4.Range("B2") '513330'
$wsQ4.Range("C2").Value = '华夏恒生互联网科技业ETF（QDII）'
Set-TextValue $wsQ4.Range("D2") '246.27'
Set-TextValue $wsQ4.Range("E2") '96.18'
Set-TextValue $wsQ4.Range("F2") '6.94'
Set-TextValue $wsQ4.Range("G2") '17.0911'
$wsQ4.Range("H2").Value = 6
$wsQ4.Range("A3").Value = 1
Set-TextValue $wsQ4.Range("B3") '513050'
$wsQ4.Range("C3").Value = '易方达中证海外中国互联网50（QDII）ETF'
Set-TextValue $wsQ4.Range("D3") '364.78'
Set-TextValue $wsQ4.Range("E3") '98.95'
Set-TextValue $wsQ4.Range("F3") '2.81'
Set-TextValue $wsQ4.Range("G3") '10.2503'
$wsQ4.Range("H3").Value = 7
$wsQ4.Range("A4").Value = 2
Set-TextValue $wsQ4.Range("B4") '513180'
$wsQ4.Range("C4").Value = '华夏恒生科技ETF（QDII）'
Set-TextValue $wsQ4.Range("D4") '168.55'
Set-TextValue $wsQ4.Range("E4") '96.27'
Set-TextValue $wsQ4.Range("F4") '4.17'
Set-TextValue $wsQ4.Range("G4") '7.0285'
$wsQ4.Range("H4").Value = 9
$wsQ4.Range("A5").Value = 3
Set-TextValue $wsQ4.Range("B5") '164906'
$wsQ4.Range("C5").Value = '交银施罗德中证海外中国互联网指数（QDII-LOF）'
Set-TextValue $wsQ4.Range("D5") '113.53'
Set-TextValue $wsQ4.Range("E5") '92.20'
Set-TextValue $wsQ4.Range("F5") '3.76'
Set-TextValue $wsQ4.Range("G5") '4.2687'
$wsQ4.Range("H5").Value = 8
$wsQ4.Range("A6").Value = 4
Set-TextValue $wsQ4.Range("B6") '159605'
$wsQ4.Range("C6").Value = '广发中证海外中国互联网30（QDII-ETF）'
Set-TextValue $wsQ4.Range("D6") '58.13'
Set-TextValue $wsQ4.Range("E6") '99.47'
Set-TextValue $wsQ4.Range("F6") '5.19'
Set-TextValue $wsQ4.Range("G6") '3.0169'
$wsQ4.Range("H6").Value = 7
$wsQ4.Range("A7").Value = 5
Set-TextValue $wsQ4.Range("B7") '513010'
$wsQ4.Range("C7").Value = '易方达恒生科技ETF（QDII）'
Set-TextValue $wsQ4.Range("D7") '39.78'
Set-TextValue $wsQ4.Range("E7") '93.48'
Set-TextValue $wsQ4.Range("F7") '4.41'
Set-TextValue $wsQ4.Range("G7") '1.7543'
$wsQ4.Range("H7").Value = 8
$wsQ4.Range("A8").Value = 6
Set-TextValue $wsQ4.Range("B8") '159607'
$wsQ4.Range("C8").Value = '嘉实中证海外中国互联网30ETF（QDII）'
Set-TextValue $wsQ4.Range("D8") '17.02'
Set-TextValue $wsQ4.Range("E8") '99.34'
Set-TextValue $wsQ4.Range("F8") '5.18'
Set-TextValue $wsQ4.Range("G8") '0.8816'
$wsQ4.Range("H8").Value = 7
$wsQ4.Range("A9").Value = 7
Set-TextValue $wsQ4.Range("B9") '159740'
$wsQ4.Range("C9").Value = '大成恒生科技ETF（QDII）'
Set-TextValue $wsQ4.Range("D9") '17.57'
Set-TextValue $wsQ4.Range("E9") '95.07'
Set-TextValue $wsQ4.Range("F9") '4.48'
Set-TextValue $wsQ4.Range("G9") '0.7871'
$wsQ4.Range("H9").Value = 8
$wsQ4.Range("A10").Value = 8
Set-TextValue $wsQ4.Range("B10") '513580'
$wsQ4.Range("C10").Value = '华安恒生科技ETF（QDII）'
Set-TextValue $wsQ4.Range("D10") '9.11'
Set-TextValue $wsQ4.Range("E10") '95.72'
Set-TextValue $wsQ4.Range("F10") '4.51'
Set-TextValue $wsQ4.Range("G10") '0.4109'
$wsQ4.Range("H10").Value = 8
$wsQ4.Range("A11").Value = 9
Set-TextValue $wsQ4.Range("B11") '012805'
$wsQ4.Range("C11").Value = '广发恒生科技指数（QDII）C'
Set-TextValue $wsQ4.Range("D11") '10.05'
Set-TextValue $wsQ4.Range("E11") '85.46'
Set-TextValue $wsQ4.Range("F11") '4.05'
Set-TextValue $wsQ4.Range("G11") '0.4070'
$wsQ4.Range("H11").Value = 8
$wsQ4.Range("A12").Value = 10
Set-TextValue $wsQ4.Range("B12") '012208'
$wsQ4.Range("C12").Value = '华夏港股前沿经济混合（QDII）A'
Set-TextValue $wsQ4.Range("D12") '10.38'
Set-TextValue $wsQ4.Range("E12") '92.91'
Set-TextValue $wsQ4.Range("F12") '3.77'
Set-TextValue $wsQ4.Range("G12") '0.3913'
$wsQ4.Range("H12").Value = 9
$wsQ4.Range("A13").Value = 11
Set-TextValue $wsQ4.Range("B13") '159742'
$wsQ4.Range("C13").Value = '博时恒生科技ETF（QDII）'
Set-TextValue $wsQ4.Range("D13") '8.45'
Set-TextValue $wsQ4.Range("E13") '94.56'
Set-TextValue $wsQ4.Range("F13") '4.46'
Set-TextValue $wsQ4.Range("G13") '0.3769'
$wsQ4.Range("H13").Value = 8
$wsQ4.Range("A14").Value = 12
Set-TextValue $wsQ4.Range("B14") '159741'
$wsQ4.Range("C14").Value = '嘉实恒生科技ETF（QDII）'
Set-TextValue $wsQ4.Range("D14") '4.64'
Set-TextValue $wsQ4.Range("E14") '99.83'
Set-TextValue $wsQ4.Range("F14") '4.70'
Set-TextValue $wsQ4.Range("G14") '0.2181'
$wsQ4.Range("H14").Value = 8
$wsQ4.Range("A15").Value = 13
Set-TextValue $wsQ4.Range("B15") '862001'
$wsQ4.Range("C15").Value = '光大阳光香港精选混合（QDII）A 人民币'
Set-TextValue $wsQ4.Range("D15") '3.90'
Set-TextValue $wsQ4.Range("E15") '92.65'
Set-TextValue $wsQ4.Range("F15") '5.57'
Set-TextValue $wsQ4.Range("G15") '0.2172'
$wsQ4.Range("H15").Value = 5
$wsQ4.Range("A16").Value = 14
Set-TextValue $wsQ4.Range("B16") '862011'
$wsQ4.Range("C16").Value = '光大阳光香港精选混合（QDII）A 美元'
Set-TextValue $wsQ4.Range("D16") '3.90'
Set-TextValue $wsQ4.Range("E16") '92.65'
Set-TextValue $wsQ4.Range("F16") '5.57'
Set-TextValue $wsQ4.Range("G16") '0.2172'
$wsQ4.Range("H16").Value = 5
$wsQ4.Range("A17").Value = 15
Set-TextValue $wsQ4.Range("B17") '862012'
$wsQ4.Range("C17").Value = '光大阳光香港精选混合（QDII）C 人民币'
Set-TextValue $wsQ4.Range("D17") '3.90'
Set-TextValue $wsQ4.Range("E17") '92.65'
Set-TextValue $wsQ4.Range("F17") '5.57'
Set-TextValue $wsQ4.Range("G17") '0.2172'
$wsQ4.Range("H17").Value = 5
$wsQ4.Range("A18").Value = 16
Set-TextValue $wsQ4.Range("B18") '513260'
$wsQ4.Range("C18").Value = '汇添富恒生科技ETF（QDII）'
Set-TextValue $wsQ4.Range("D18") '4.48'
Set-TextValue $wsQ4.Range("E18") '98.09'
Set-TextValue $wsQ4.Range("F18") '4.63'
Set-TextValue $wsQ4.Range("G18") '0.2074'
$wsQ4.Range("H18").Value = 8
$wsQ4.Range("A19").Value = 17
Set-TextValue $wsQ4.Range("B19") '012804'
$wsQ4.Range("C19").Value = '广发恒生科技指数（QDII）A'
Set-TextValue $wsQ4.Range("D19") '4.88'
Set-TextValue $wsQ4.Range("E19") '85.46'
Set-TextValue $wsQ4.Range("F19") '4.05'
Set-TextValue $wsQ4.Range("G19") '0.1976'
$wsQ4.Range("H19").Value = 8
$wsQ4.Range("A20").Value = 18
Set-TextValue $wsQ4.Range("B20") '159747'
$wsQ4.Range("C20").Value = '南方中证香港科技ETF（QDII）'
Set-TextValue $wsQ4.Range("D20") '3.55'
Set-TextValue $wsQ4.Range("E20") '99.72'
Set-TextValue $wsQ4.Range("F20") '4.84'
Set-TextValue $wsQ4.Range("G20") '0.1718'
$wsQ4.Range("H20").Value = 7
$wsQ4.Range("A21").Value = 19
Set-TextValue $wsQ4.Range("B21") '012379'
$wsQ4.Range("C21").Value = '创金合信港股互联网3个月持有期混合（QDII）A'
Set-TextValue $wsQ4.Range("D21") '3.24'
Set-TextValue $wsQ4.Range("E21") '88.37'
Set-TextValue $wsQ4.Range("F21") '3.99'
Set-TextValue $wsQ4.Range("G21") '0.1293'
$wsQ4.Range("H21").Value = 7
$wsQ4.Range("A22").Value = 20
Set-TextValue $wsQ4.Range("B22") '513890'
$wsQ4.Range("C22").Value = '上投摩根恒生科技ETF（QDII）'
Set-TextValue $wsQ4.Range("D22") '2.28'
Set-TextValue $wsQ4.Range("E22") '97.27'
Set-TextValue $wsQ4.Range("F22") '4.61'
Set-TextValue $wsQ4.Range("G22") '0.1051'
$wsQ4.Range("H22").Value = 10
$wsQ4.Range("A23").Value = 21
Set-TextValue $wsQ4.Range("B23") '012380'
$wsQ4.Range("C23").Value = '创金合信港股互联网3个月持有期混合（QDII）C'
Set-TextValue $wsQ4.Range("D23") '1.18'
Set-TextValue $wsQ4.Range("E23") '88.37'
Set-TextValue $wsQ4.Range("F23") '3.99'
Set-TextValue $wsQ4.Range("G23") '0.0471'
$wsQ4.Range("H23").Value = 7
$wsQ4.Range("A24").Value = 22
Set-TextValue $wsQ4.Range("B24") '513220'
$wsQ4.Range("C24").Value = '招商中证全球中国互联网ETF（QDII）'
Set-TextValue $wsQ4.Range("D24") '1.10'
Set-TextValue $wsQ4.Range("E24") '98.97'
Set-TextValue $wsQ4.Range("F24") '3.99'
Set-TextValue $wsQ4.Range("G24") '0.0439'
$wsQ4.Range("H24").Value = 8
$wsQ4.Range("A25").Value = 23
Set-TextValue $wsQ4.Range("B25") '159750'
$wsQ4.Range("C25").Value = '招商中证香港科技ETF（QDII）'
Set-TextValue $wsQ4.Range("D25") '1.06'
Set-TextValue $wsQ4.Range("E25") '98.79'
Set-TextValue $wsQ4.Range("F25") '4.12'
Set-TextValue $wsQ4.Range("G25") '0.0437'
$wsQ4.Range("H25").Value = 9
$wsQ4.Range("A26").Value = 24
Set-TextValue $wsQ4.Range("B26") '513380'
$wsQ4.Range("C26").Value = '广发恒生科技（QDII-ETF）'
Set-TextValue $wsQ4.Range("D26") '1.00'
Set-TextValue $wsQ4.Range("E26") '91.23'
Set-TextValue $wsQ4.Range("F26") '4.33'
Set-TextValue $wsQ4.Range("G26") '0.0433'
$wsQ4.Range("H26").Value = 8
$wsQ4.Range("A27").Value = 25
Set-TextValue $wsQ4.Range("B27") '012209'
$wsQ4.Range("C27").Value = '华夏港股前沿经济混合（QDII）C'
Set-TextValue $wsQ4.Range("D27") '0.38'
Set-TextValue $wsQ4.Range("E27") '92.91'
Set-TextValue $wsQ4.Range("F27") '3.77'
Set-TextValue $wsQ4.Range("G27") '0.0143'
$wsQ4.Range("H27").Value = 9
$wsQ4.Range("A28").Value = 26
Set-TextValue $wsQ4.Range("B28") '378006'
$wsQ4.Range("C28").Value = '上投摩根全球新兴市场混合（QDII）'
Set-TextValue $wsQ4.Range("D28") '0.44'
Set-TextValue $wsQ4.Range("E28") '86.49'
Set-TextValue $wsQ4.Range("F28") '1.59'
Set-TextValue $wsQ4.Range("G28") '0.0070'
$wsQ4.Range("H28").Value = 9
Write-Host "Workbook now has" $wb.Worksheets.Count "sheets"
foreach ($s in $wb.Worksheets) {
    Write-Host " -" $s.Name
}
